$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99-115 down to 100-116
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with the new weekly price record
$ws.Range("A99").Value = 10
$ws.Range("B99").Value = "Vega Modelo de Temuco"
$ws.Range("C99").Value = "La Araucanía"
$ws.Range("D99").Value = 45218
$ws.Range("E99").Value = 9
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100107
$ws.Range("H99").Value = "Otros"
$ws.Range("I99").Value = 100107011
$ws.Range("J99").Value = "Tuna"
$ws.Range("K99").Value = "Sin especificar"
$ws.Range("L99").Value = "Especial"
$ws.Range("M99").Value = 140
$ws.Range("N99").Value = 36000
$ws.Range("O99").Value = 36000
$ws.Range("P99").Value = 36000
$ws.Range("Q99").Value = "$/caja 16 kilos"
$ws.Range("R99").Value = "Provincia de Los Andes"
$ws.Range("S99").Value = 2250
$ws.Range("T99").Value = 16
